# Auto-generated Excel COM-interop script
# Updates the cryptos price list (columns D = Price, E = Volume(1h)) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.143.61"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("D3").Value = "2.439.44"
$ws.Range("E3").Value = "  +4.61%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.32"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.64"
$ws.Range("E6").Value = "  +5.54%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "2.435.35"
$ws.Range("E9").Value = "  +4.57%  "
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").Value = "  +4.22%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("E13").Value = "  +4.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.83"
$ws.Range("E14").Value = "  +9.25%  "
$ws.Range("D15").Value = "2.867.09"
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").Value = "62.038.16"
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("E17").Value = "  +5.54%  "
$ws.Range("D18").Value = "2.435.77"
$ws.Range("E18").Value = "  +4.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.14"
$ws.Range("E19").Value = "  +5.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.08"
$ws.Range("E20").Value = "  +9.84%  "
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.85"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.15"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.174"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +11.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.29"
$ws.Range("E28").Value = "  +5.84%  "
$ws.Range("E29").Value = "  +12.76%  "
$ws.Range("E30").Value = "  +7.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.80"
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.34"
$ws.Range("E32").Value = "  +6.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "171.12"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.45"
$ws.Range("E34").Value = "  +5.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.397"
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "374.92"
$ws.Range("E36").Value = "  +16.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.59"
$ws.Range("E37").Value = "  +4.18%  "
$ws.Range("E38").Value = "  +10.17%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.69"
$ws.Range("E41").Value = "  +10.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.15"
$ws.Range("E42").Value = "  +3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.70"
$ws.Range("E43").Value = "  +6.86%  "
$ws.Range("E44").Value = "  +5.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.70"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0958"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.590"
$ws.Range("E47").Value = "  +4.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0521"
$ws.Range("E48").Value = "  +5.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.00"
$ws.Range("E49").Value = "  +6.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0222"
$ws.Range("E50").Value = "  +3.92%  "
$ws.Range("E51").Value = "  +11.52%  "
